$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the "Materias primas" ingredients for "pie de limon" (row 6):
# move "1.0-crema," from the 3rd position to the end.
$ws.Range("C6").Value = "5.0-merengue,2.0-limon,5.0-huevos,4.0-harinita,1.0-crema,"

# Add new product row 8
$ws.Range("A8").Value = "m"
$ws.Range("B8").Value = 2.0
$ws.Range("C8").Value = "2.0-harinita,"
$ws.Range("D8").Value = -3.0
